$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.814.66"
$ws.Range("E2").Value = "  +5.47%  "
$ws.Range("D3").Value = "'2.750.15"
$ws.Range("E3").Value = "  +4.38%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'581.34"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'153.38"
$ws.Range("E6").Value = "  +6.05%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").Value = "'2.751.21"
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("D10").Value = "'6.73"
$ws.Range("E10").Value = "  +3.25%  "
$ws.Range("D13").Value = "'0.161"
$ws.Range("E13").Value = "  +3.90%  "
$ws.Range("D14").Value = "'3.208.63"
$ws.Range("E14").Value = "  +3.82%  "
$ws.Range("D15").Value = "'26.31"
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'63.403.99"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("E17").Value = "  +7.17%  "
$ws.Range("D18").Value = "'2.743.83"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("E20").Value = "  +3.38%  "
$ws.Range("D21").Value = "'360.42"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").Value = "'6.99"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "'0.538"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "'65.88"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  +4.78%  "
$ws.Range("D27").Value = "'8.58"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'0.0₃0896"
$ws.Range("E29").Value = "  +12.43%  "
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("E31").Value = "  +6.63%  "
$ws.Range("D32").Value = "'172.97"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("E33").Value = "  +14.35%  "
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").Value = "  +4.64%  "
$ws.Range("E36").Value = "  +7.73%  "
$ws.Range("E37").Value = "  +8.60%  "
$ws.Range("E38").Value = "  +9.35%  "
$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "  +15.53%  "
$ws.Range("D40").Value = "'345.43"
$ws.Range("E40").Value = "  +4.61%  "
$ws.Range("D41").Value = "'4.21"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").Value = "'39.05"
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  +8.04%  "
$ws.Range("D44").Value = "'21.77"
$ws.Range("E44").Value = "  +8.46%  "
$ws.Range("D45").Value = "'21.86"
$ws.Range("E45").Value = "  +6.19%  "
$ws.Range("E46").Value = "  +6.17%  "
$ws.Range("D47").Value = "'139.13"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("E48").Value = "  +5.49%  "
$ws.Range("D49").Value = "'0.0256"
$ws.Range("E49").Value = "  +4.78%  "
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  +0.01%  "
